$d = $word.ActiveDocument

$replacements = @(
    @{old = "713÷3="; new = "303÷3="},
    @{old = "574÷3="; new = "946÷4="},
    @{old = "136÷2="; new = "988÷3="},
    @{old = "374÷4="; new = "773÷4="},
    @{old = "455÷5="; new = "520÷3="},
    @{old = "699÷9="; new = "981÷2="},
    @{old = "179÷7="; new = "626÷5="},
    @{old = "154÷6="; new = "800÷2="},
    @{old = "867÷7="; new = "121÷7="},
    @{old = "853÷6="; new = "467÷2="},
    @{old = "532÷6="; new = "290÷7="},
    @{old = "432÷8="; new = "103÷7="},
    @{old = "958÷6="; new = "812÷3="},
    @{old = "736÷6="; new = "812÷9="},
    @{old = "995÷6="; new = "375÷5="},
    @{old = "606÷2="; new = "178÷4="},
    @{old = "722÷4="; new = "808÷4="},
    @{old = "335÷6="; new = "444÷9="},
    @{old = "135÷9="; new = "118÷5="},
    @{old = "717÷6="; new = "575÷9="},
    @{old = "173÷9="; new = "439÷5="},
    @{old = "898÷6="; new = "401÷8="},
    @{old = "935÷7="; new = "900÷6="},
    @{old = "619÷6="; new = "997÷8="},
    @{old = "347÷5="; new = "989÷5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
